$d = $word.ActiveDocument

function Isolate-And-Split($foundRange, $newText, $splitOffsets) {
    # Replace the text of $foundRange with $newText, then force it to become
    # its own run (isolated from whatever run precedes it in the paragraph),
    # and further force additional run boundaries at each offset in
    # $splitOffsets (character offsets from the start of $newText). The
    # isolation/boundary trick toggles Bold on/off (net no-op formatting
    # change) on a sub-range, which is enough to stop the engine from
    # merging two textually-adjacent same-formatted runs back together.
    $start = $foundRange.Start
    $foundRange.Text = $newText
    $len = $newText.Length
    $whole = $d.Range($start, $start + $len)
    $whole.Bold = 1
    $whole.Bold = 0
    foreach ($off in $splitOffsets) {
        $sub = $d.Range($start + $off, $start + $len)
        $sub.Bold = 1
        $sub.Bold = 0
    }
}

# ---------------------------------------------------------------------
# Paragraph 23, piece 1: "进行阐释，画图请使用"
#   -> "进行阐释，画图" | "建议" | "使用"
# ---------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("进行阐释，画图请使用", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Isolate-And-Split $r1 "进行阐释，画图建议使用" @(7, 9)

# ---------------------------------------------------------------------
# Paragraph 23, piece 2: "。内容格式可以按照附件中的"
#   -> "。" | "matplotlib" | "也可以" | "，" | "内容格式可以按照附件中的"
# ---------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("。内容格式可以按照附件中的", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Isolate-And-Split $r2 "。matplotlib也可以，内容格式可以按照附件中的" @(1, 11, 14, 15)

# ---------------------------------------------------------------------
# Paragraph 24, piece 1: "请登记在" -> "可以" | "登记在"
# ---------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("请登记在", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Isolate-And-Split $r3 "可以登记在" @(2)

# ---------------------------------------------------------------------
# Paragraph 24, piece 2: "供大家分享纠错。" -> "供大家分享" | "学习" | "。"
# ---------------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("供大家分享纠错。", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Isolate-And-Split $r4 "供大家分享学习。" @(5, 7)
